$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "63.424.23"
$ws.Range("E2").Value = "  -1.34%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.065.97"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.24%  "

# Row 5 - BNB
Set-TextValue "D5" "588.58"
$ws.Range("E5").Value = "  -0.65%  "

# Row 6 - Solana
Set-TextValue "D6" "154.76"
$ws.Range("E6").Value = "  +4.75%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.07%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.80%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "3.064.03"
$ws.Range("E9").Value = "  -2.49%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -4.17%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -1.46%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.450"
$ws.Range("E12").Value = "  -1.76%  "

# Row 13 - Avalanche
Set-TextValue "D13" "36.93"
$ws.Range("E13").Value = "  -1.48%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -4.09%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  -2.07%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue "D16" "3.574.56"
$ws.Range("E16").Value = "  -2.76%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "63.478.82"
$ws.Range("E17").Value = "  -0.90%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  -2.09%  "

# Row 19
Set-TextValue "D19" "3.063.15"
$ws.Range("E19").Value = "  -2.75%  "

# Row 20
Set-TextValue "D20" "472.95"
$ws.Range("E20").Value = "  +0.96%  "

# Row 21
Set-TextValue "D21" "14.31"
$ws.Range("E21").Value = "  -1.19%  "

# Row 22
Set-TextValue "D22" "0.705"
$ws.Range("E22").Value = "  -3.94%  "

# Row 23
Set-TextValue "D23" "7.52"
$ws.Range("E23").Value = "  -1.30%  "

# Row 24
Set-TextValue "D24" "2.42"
$ws.Range("E24").Value = "  +0.14%  "

# Row 25
Set-TextValue "D25" "80.70"
$ws.Range("E25").Value = "  -0.74%  "

# Row 26
Set-TextValue "D26" "12.82"
$ws.Range("E26").Value = "  -2.63%  "

# Row 27
Set-TextValue "D27" "10.41"
$ws.Range("E27").Value = "  +4.09%  "

# Row 28
Set-TextValue "D28" "0.999"
$ws.Range("E28").Value = "  -0.23%  "

# Row 29
Set-TextValue "D29" "7.49"
$ws.Range("E29").Value = "  +2.83%  "

# Row 30
$ws.Range("E30").Value = "  -0.26%  "

# Row 31
$ws.Range("E31").Value = "  -2.30%  "

# Row 32
$ws.Range("E32").Value = "  -3.14%  "

# Row 33
$ws.Range("E33").Value = "  -2.16%  "

# Row 34
Set-TextValue "D34" "27.11"
$ws.Range("E34").Value = "  -2.61%  "

# Row 35 - PEPE
Set-TextValue "D35" "0.0₃0819"
$ws.Range("E35").Value = "  -4.97%  "

# Row 36
$ws.Range("E36").Value = "  -2.31%  "

# Row 37
Set-TextValue "D37" "3.32"
$ws.Range("E37").Value = "  +1.10%  "

# Row 38
Set-TextValue "D38" "5.99"
$ws.Range("E38").Value = "  -3.09%  "

# Row 39
$ws.Range("E39").Value = "  -3.96%  "

# Row 40 - was OKB, now Cosmos
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D40" "9.25"
$ws.Range("E40").Value = "  -0.55%  "

# Row 41 - was Cosmos, now OKB
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D41" "50.51"
$ws.Range("E41").Value = "  -1.52%  "

# Row 42 - Bittensor
Set-TextValue "D42" "437.28"
$ws.Range("E42").Value = "  -5.55%  "

# Row 43 - TheGraph
$ws.Range("E43").Value = "  -2.03%  "

# Row 44 - Arweave
Set-TextValue "D44" "40.65"
$ws.Range("E44").Value = "  +1.79%  "

# Row 45 - Kaspa
Set-TextValue "D45" "0.112"
$ws.Range("E45").Value = "  +3.41%  "

# Row 46 - VeChain
Set-TextValue "D46" "0.0359"
$ws.Range("E46").Value = "  -4.18%  "

# Row 47 - Maker
Set-TextValue "D47" "2.793.91"
$ws.Range("E47").Value = "  -3.42%  "

# Row 48 - Monero
Set-TextValue "D48" "130.01"
$ws.Range("E48").Value = "  -2.36%  "

# Row 49 - USDe
$ws.Range("E49").Value = "  +0.03%  "

# Row 50 - InjectiveProtocol
Set-TextValue "D50" "25.14"
$ws.Range("E50").Value = "  +4.03%  "

# Row 51 - ThetaToken
$ws.Range("E51").Value = "  -0.18%  "
